$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数")
$updates = @{
    2  = 285
    4  = 144
    6  = 481
    7  = 1366
    8  = 549
    9  = 97
    10 = 172
    11 = 119
    12 = 173
    13 = 97
    14 = 146
    15 = 136
}

# Both "展览" and "全部类型" sheets contain identical data that need updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
